$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.5
$ws.Range("H2").Value = 3.3
$ws.Range("I2").Value = 2.8
$ws.Range("N2").Value = 2
$ws.Range("O2").Value = 1.9
$ws.Range("J3").Value = 1.03
$ws.Range("K3").Value = 17
$ws.Range("L3").Value = 1.17
$ws.Range("M3").Value = 5
$ws.Range("G13").Value = 1.55
$ws.Range("I13").Value = 5
$ws.Range("N13").Value = 1.53
$ws.Range("O13").Value = 2.4
$ws.Range("U13").Value = 9
$ws.Range("W13").Value = 12
$ws.Range("AF13").Value = 29
$ws.Range("G14").Value = 2.45
$ws.Range("I14").Value = 3.1
$ws.Range("T14").Value = 7
$ws.Range("U14").Value = 11
$ws.Range("W14").Value = 23
$ws.Range("AD14").Value = 1250
$ws.Range("H15").Value = 3.65
$ws.Range("N15").Value = 1.83
$ws.Range("O15").Value = 1.78
$ws.Range("T15").Value = 5.6
$ws.Range("V15").Value = 6.9
$ws.Range("W15").Value = 9.5
$ws.Range("Z15").Value = 9.5
$ws.Range("AF15").Value = 23
$ws.Range("AH15").Value = 70
$ws.Range("G17").Value = 1.6
$ws.Range("H17").Value = 3.5
$ws.Range("I17").Value = 5.1
$ws.Range("T17").Value = 5.4
$ws.Range("U17").Value = 6.1
$ws.Range("V17").Value = 6.8
$ws.Range("W17").Value = 9.75
$ws.Range("X17").Value = 10.75
$ws.Range("Y17").Value = 22
$ws.Range("Z17").Value = 9.5
$ws.Range("AA17").Value = 6
$ws.Range("AB17").Value = 13.5
$ws.Range("AC17").Value = 60
$ws.Range("AD17").Value = 450
$ws.Range("AE17").Value = 11
$ws.Range("AF17").Value = 24
$ws.Range("AG17").Value = 13.5
$ws.Range("AH17").Value = 75
$ws.Range("AI17").Value = 45
$ws.Range("AJ17").Value = 45
$ws.Range("H18").Value = 4.1
$ws.Range("I18").Value = 6
$ws.Range("T18").Value = 6.1
$ws.Range("V18").Value = 7
$ws.Range("Z18").Value = 12
$ws.Range("AA18").Value = 7.2
$ws.Range("AE18").Value = 14
$ws.Range("AF18").Value = 30
$ws.Range("AH18").Value = 90
$ws.Range("G19").Value = 1.09
$ws.Range("H19").Value = 7.7
$ws.Range("I19").Value = 17
$ws.Range("T19").Value = 9
$ws.Range("U19").Value = 5.9
$ws.Range("V19").Value = 10.75
$ws.Range("W19").Value = 5.3
$ws.Range("X19").Value = 9.75
$ws.Range("Y19").Value = 32
$ws.Range("Z19").Value = 21
$ws.Range("AA19").Value = 16.5
$ws.Range("AB19").Value = 35
$ws.Range("AC19").Value = 150
$ws.Range("AE19").Value = 45
$ws.Range("AF19").Value = 150
$ws.Range("AG19").Value = 55
$ws.Range("AI19").Value = 250
$ws.Range("AJ19").Value = 150
$ws.Range("N21").Value = 2.3
$ws.Range("O21").Value = 1.6
$ws.Range("G24").Value = 1.22
$ws.Range("H24").Value = 6
$ws.Range("I24").Value = 9.5
$ws.Range("J24").Value = 19
$ws.Range("K24").Value = 1.03
$ws.Range("L24").Value = 1.13
$ws.Range("M24").Value = 5.5
$ws.Range("N24").Value = 1.44
$ws.Range("O24").Value = 2.63
$ws.Range("P24").Value = 1.22
$ws.Range("Q24").Value = 4
$ws.Range("R24").Value = 1.91
$ws.Range("S24").Value = 1.8
$ws.Range("T24").Value = 9.5
$ws.Range("U24").Value = 7.5
$ws.Range("Z24").Value = 19
$ws.Range("AC24").Value = 51
$ws.Range("AD24").Value = 600
$ws.Range("AE24").Value = 29
$ws.Range("AF24").Value = 51
$ws.Range("AG24").Value = 29
$ws.Range("AH24").Value = 126
$ws.Range("AI24").Value = 67
$ws.Range("G25").Value = 1.95
$ws.Range("H25").Value = 3.6
$ws.Range("I25").Value = 3.4
$ws.Range("R25").Value = 1.73
$ws.Range("S25").Value = 2
$ws.Range("U25").Value = 10
$ws.Range("V25").Value = 9
$ws.Range("W25").Value = 17
$ws.Range("X25").Value = 15
$ws.Range("Y25").Value = 23
$ws.Range("AB25").Value = 15
$ws.Range("AE25").Value = 12
$ws.Range("AF25").Value = 19
$ws.Range("AG25").Value = 13
$ws.Range("AH25").Value = 41
$ws.Range("AI25").Value = 26
$ws.Range("AJ25").Value = 34
$ws.Range("M30").Value = 5.1
$ws.Range("R30").Value = 1.82
$ws.Range("T30").Value = 9.25
$ws.Range("U30").Value = 7.2
$ws.Range("W30").Value = 8
$ws.Range("AA30").Value = 11.5
$ws.Range("AF30").Value = 100
$ws.Range("AJ30").Value = 75
